$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.24806121495476674
$ws.Range("A2").Value = -0.0059999999592719178
$ws.Range("A3").Value = -0.0039999999639857009
$ws.Range("A4").Value = -0.0079999999346185291
$ws.Range("A5").Value = -0.0029999999630136998
$ws.Range("A6").Value = -0.0019999999601285623
$ws.Range("A7").Value = -0.0099999999084672275
$ws.Range("A8").Value = -0.0099999999079027901
$ws.Range("A9").Value = -0.00199999995929101
$ws.Range("A10").Value = -0.0019999999597075657
$ws.Range("A11").Value = -0.0029999999533174559
$ws.Range("A12").Value = -0.0034999999504616852
$ws.Range("A13").Value = 0.027869823833918339
$ws.Range("A14").Value = -0.0079999999229052321
$ws.Range("A15").Value = -0.00099999996724786655
$ws.Range("A16").Value = 0.00011434601366744346
$ws.Range("A17").Value = -0.0019999999594944029
$ws.Range("A18").Value = -0.0039999999464104263
$ws.Range("A19").Value = -0.049311834351340345
$ws.Range("A20").Value = -0.0039999999622963855
$ws.Range("A21").Value = -0.00399999996040723
$ws.Range("A22").Value = -0.035032666332921458
$ws.Range("A23").Value = -0.0049999999555998542
$ws.Range("A24").Value = -0.019999999855594197
$ws.Range("A25").Value = -0.019999999853611783
$ws.Range("A26").Value = 0.029175098240012076
$ws.Range("A27").Value = -0.0024999999526396088
$ws.Range("A28").Value = -0.001999999942773556
$ws.Range("A29").Value = -0.0069999999010645908
$ws.Range("A30").Value = -0.059999999556914485
$ws.Range("A31").Value = -0.0069999998932761542
$ws.Range("A32").Value = -0.0099999998735569307
$ws.Range("A33").Value = -0.0039999999118336405
